$p = $ppt.ActivePresentation

# --- Slide 10 ("Class Diagram") -------------------------------------------
# The title placeholder on the last existing slide was empty; give it a title.
$slide10 = $p.Slides.Item(10)
$title10 = $slide10.Shapes.Item(1)
$title10.TextFrame.TextRange.Text = "Class Diagram"
$title10.TextFrame.TextRange.LanguageID = "en-GB"

# --- New slide 11 ("Live Demo") --------------------------------------------
# Add a new slide using the "Section Header" layout (same one used by the
# "System Design" / "Design Patterns" slides earlier in the deck), which
# provides a Title placeholder plus a body/text placeholder (idx 1).
$master = $p.SlideMaster
$sectionHeaderLayout = $master.CustomLayouts.Item(3)
$slide11 = $p.Slides.AddSlide($p.Slides.Count + 1, $sectionHeaderLayout)

$title11 = $slide11.Shapes.Item(1)
$title11.Name = "Title 6"
$title11.TextFrame.TextRange.Text = "Live Demo"
$title11.TextFrame.TextRange.LanguageID = "en-GB"

$body11 = $slide11.Shapes.Item(2)
$body11.Name = "Text Placeholder 7"
$body11.TextFrame.TextRange.Text = "Java desktop app & HTML web app"
$body11.TextFrame.TextRange.LanguageID = "en-GB"
